# edit.ps1 — applies the "Agregué la presentación personal" commit:
#   1) Fix the spell-check artifact around "mugs" by collapsing the three
#      runs (and the two <w:proofErr/> markers around "mugs") that spell
#      out "la cual vende ... entre otros." into a single plain run.
#   2) Append Milena Castaño's personal-presentation block (blank list
#      paragraph, bulleted name, age, bio, role) right after the
#      "Rol: Analista" paragraph that currently ends the document.

$d = $word.ActiveDocument

# --- Part 1: merge the "la cual vende ... entre otros." runs -------------
$sentence = "la cual vende productos en madera tales como libretas, llaveros, " + `
            "lapiceros y además mugs personalizados, entre otros."

$seek = $d.Content.Duplicate
$found = $seek.Find.Execute($sentence, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($seek.Start, $seek.End)
    $mergedXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        '<w:p><w:r><w:t>' + $sentence + '</w:t></w:r></w:p>' + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($mergedXml)
}

# --- Part 2: append Milena Castaño's presentation block -------------------
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Rol: Analista") | Out-Null
$insertPoint = $d.Range($anchor.End, $anchor.End)

$newParasXml = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="360"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Milena Castaño</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="360"/></w:pPr><w:r><w:t>Tengo 26 años</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="360"/></w:pPr><w:r><w:t>Soy Ingeniera de Sistemas e Informática de la Universidad Nacional de Colombia y Especialista en Analítica de la misma universidad.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="360"/></w:pPr><w:r><w:t>Rol: Administradora de Base de datos</w:t></w:r></w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    $newParasXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($packageXml)

Write-Output "Applied presentation edits."
